# "Now using cython as default filter" - update example dates used by the
# filter/import examples on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# E12 (valid_from for name-1): shift from 2019-01-01 to 2018-01-01
$ws.Range("E12").Value = Get-Date -Year 2018 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0

# F16 (valid_until for name-4): add 2018-12-31
$ws.Range("F16").Value = Get-Date -Year 2018 -Month 12 -Day 31 -Hour 0 -Minute 0 -Second 0

# E17 / F17 (valid_from / valid_until for name-4, row 17): add 2019-01-01 / 2019-12-31
$ws.Range("E17").Value = Get-Date -Year 2019 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0
$ws.Range("F17").Value = Get-Date -Year 2019 -Month 12 -Day 31 -Hour 0 -Minute 0 -Second 0

# Leave the UI selection on E12, matching the saved view state
$ws.Range("E12").Select()
